$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.711.37"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "1.646.67"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "213.06"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "23.27"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.879.69"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.643.55"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "64.87"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "27.707.59"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "231.95"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "7.65"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +10.95%  "
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").Value = "150.04"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").Value = "6.93"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Value = "15.66"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Value = "1.442.69"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "0.0166"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("E40").Value = "  +12.12%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "67.63"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").Value = "1.789.05"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "85.66"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0989"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.79"
$ws.Range("E51").Value = "  +1.37%  "
